$d = $word.ActiveDocument

# --- paragraph 6 (0-based body index 5) ---
$p = $d.Paragraphs.Item(6)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("1.Conhecer modelos de elaboração de projetos: Conceitos gerais, diferentes modelos de projetos. `v2.Elaborar projetos: Definição da problemática, justificativas, objetivos e hipóteses, bases teóricas fundamentais, metodologia, cronograma, resultados esperados. `v3.Executar as etapas do projeto, buscando eventuais mudanças de direcionamento. `v4.Finalizar o projeto redigindo e avaliando os resultados finais. Apresentação de protótipo.`v5.Estudos preliminares de ampliação de escala do projeto.")

# --- paragraph 8 (0-based body index 7) ---
$p = $d.Paragraphs.Item(8)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("Disciplina integradora que visa desenvolver projetos na área de Engenharia Química, com especificidade em Processos Químicos.`vIntegrar, através de atividades de projeto contextualizado, os conhecimentos desenvolvidos nas unidades curriculares das disciplinas de Química Geral, Química Geral Experimental, Introdução à Engenharia Química e Balanço de Massa e Energia. Desenvolver competências de trabalho em equipe, comunicação oral e escrita, resolução de problemas, pensamento crítico, pensamento criativo, metodologia de desenvolvimento de projetos visando ao desenvolvimento das competências adquiridas no curso através de aplicação em projetos na área de Processos Químicos.`vAvaliação de Projeto: `v-Apresentações orais (pré-projeto, relatório preliminar, relatório final, ampliação de escala);`v-Trabalhos escritos (relatório preliminar e relatório final);`v-Avaliação pelos pares.`vMédia Final = Nota de Projeto `vMédia final mínima de aprovação = 5,0")

# --- paragraph 10 (0-based body index 9) ---
$p = $d.Paragraphs.Item(10)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("(Prova escrita + Média Final)/2         Nota Final Mínima para Aprovação= 5,0")

# --- paragraph 12 (0-based body index 11) ---
$p = $d.Paragraphs.Item(12)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("Powell, P. C., & Weenk, W. (2003). Project-Led Engineering Education. Utrecht: Lemma.`vUNESCO (2010). Engineering: Issues, Challenges and Opportunities for Development. Paris, France, United Nations Educational, Scientific and Cultural Organization. Retrieved from http://unesdoc.unesco.org/images/0018/001897/189753e.pdf`vLima, R. M., Carvalho, D., Sousa, R. M., Alves, A., Moreira, F., Mesquita, D., & Fernandes, S. (2011). Estrutura de Gestão para Planejamento e Execução de Projetos Interdisciplinares de Aprendizagem em Engenharia. In L. C. d. Campos, E. A. T. Dirani & A. L. Manrique (Eds.),Educação em Engenharia: Novas Abordagens (pp. 87-121). São Paulo, Brasil: EDUC  Editora da PUC-SP.`vFernandes, S., Flores, M. A., & Lima, R. M. (2011). A Avaliação dos Alunos no Contexto de um Projeto Interdisciplinar. In L. C. d. Campos, E. A. T. Dirani & A. L. Manrique (Eds.), Educação em Engenharia: Novas Abordagens (pp. 219-280). São Paulo, Brasil: EDUC  Editora da PUC-SP.`vATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006.`vBRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.`vBROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.`vCHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010. `vRUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill`vMIHELCIC, J. R. Fundamentals of Environmental Engineering. John Wiley & Sons, 1998.`vHIMMELBLAU, D. M. Engenharia Química Princípios e Cálculos. Prentice-Hall do Brasil. 1998.`vFELDER, R.M. & ROUSSEAU, R.W. Princípios Elementares dos Processos Químicos. LTC, 2005.`vSHREVE, R. Norris; BRINK JR, Joseph A. Indústria de processos químicos. Rio de Janeiro: Guanabara Dois, 1980.`vBRASIL, Nilo Indio do. Introdução a engenharia química. Rio de Janeiro: Interciencia/Petrobras, 2004.`vCREMASCO, Marco Aurélio. Engenharia química. Ed. Edgard Blucher, 2005.")

# --- paragraph 14 (0-based body index 13) ---
$p = $d.Paragraphs.Item(14)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("Método: 198273 - Domingos Savio Giordani`vCritério: 5817045 - Elisângela de Jesus Cândido Moraes`vNorma de recuperação: 5817344 - Livia Melo Carneiro")
$d.Range($pStart + 0, $pStart + 8).Font.Bold = 1
$d.Range($pStart + 41, $pStart + 51).Font.Bold = 1
$d.Range($pStart + 96, $pStart + 118).Font.Bold = 1

# --- paragraph 16 (0-based body index 15) ---
$p = $d.Paragraphs.Item(16)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$pStart = $full.Start
$full.Delete()
$full.InsertAfter("787307 - Luis Fernando Figueiredo Faria")
